# Hotfix: Wed Oct 16 17:38:15 RTZ 2024
# Appends new rows to the "Python" and "Bash" sheets of the database_tables workbook.
$wb = $excel.ActiveWorkbook

# ---- Sheet: Python (new rows 33-37) ----
$wsPy = $wb.Worksheets.Item("Python")

$wsPy.Range("A33").Value = 2094
$wsPy.Range("B33").Value = @'
Скрипт удаления из таблицы баз данных определенных id 
'@
$wsPy.Range("C33").Value = @'
Скрипт удаления из таблицы баз данных определенных id 
'@
$wsPy.Range("D33").Value = @'
# Подключаем библиотеку sqlite3
import sqlite3
# Подключаемся в базе данных
con = sqlite3.connect("../database1.db")
cur = con.cursor()
res = cur.execute("DELETE FROM links WHERE id IN (12,13,14)")
res.fetchall()
con.commit()
'@

$wsPy.Range("A34").Value = 2095
$wsPy.Range("B34").Value = @'
Скрипт очистки таблицы базы данных
'@
$wsPy.Range("C34").Value = @'
Скрипт очистки таблицы базы данных
'@
$wsPy.Range("D34").Value = @'
import sqlite3
con = sqlite3.connect("/database1.db")
cur = con.cursor()
res = cur.execute("DROP TABLE table_name")
con.commit()
'@

$wsPy.Range("A35").Value = 2096
$wsPy.Range("B35").Value = @'
Скрипт удаления таблицы
'@
$wsPy.Range("C35").Value = @'
Скрипт удаления таблицы
'@
$wsPy.Range("D35").Value = @'
import sqlite3
con = sqlite3.connect("/database1.db")
cur = con.cursor()
res = cur.execute("DROP TABLE table_name")
con.commit()
'@

$wsPy.Range("A36").Value = 2097
$wsPy.Range("B36").Value = @'
Скрипт вставки значений в таблицу
'@
$wsPy.Range("C36").Value = @'
Скрипт вставки значений в таблицу
'@
$wsPy.Range("D36").Value = @'
# Подключаем библиотеку sqlite3
import sqlite3
# Подключаемся в базе данных
con = sqlite3.connect("../database1.db")
cur = con.cursor()
data = (
    {"id": None, "name": "test", "link": "test"},
    {"id": None, "name": "test", "link": "test"},
    {"id": None, "name": "test", "link": "test"},
)
cur.executemany("INSERT INTO links VALUES(:id,:name, :link)", data)
con.commit()

'@

$wsPy.Range("A37").Value = 2098
$wsPy.Range("B37").Value = @'
Экземпляр приложения Flask
'@
$wsPy.Range("C37").Value = @'
Экземпляр приложения Flask с подключением внешней базы данных и пагинацией
'@
$wsPy.Range("D37").Value = @'
import flask
import pymysql.cursors
from flask_paginate import Pagination, get_page_args
app = flask.Flask(__name__)
app.secret_key = "secret key"
@app.errorhandler(404)
def page_not_found(e):
    return flask.render_template('404.html'), 404
def get_db_connection():
    conn = pymysql.connect(host='localhost',
                           port=3307,
                           user='root',
                           password='1',
                           database='test_base',
                           charset='utf8',
                           cursorclass=pymysql.cursors.DictCursor)
    return conn
def close_db_connection(conn):
    conn.close()
@app.route("/")
def index():
    conn = get_db_connection()
    with conn.cursor() as cur:
        cur.execute("SELECT * FROM test_table")
        test_list_posts = cur.fetchall()
    conn.close()
    page, per_page, offset = get_page_args(page_parameter='page',
                                           per_page_parameter='per_page')
    total = len(test_list_posts)
    def get_test_list_posts(offset=0, per_page=5):
        return test_list_posts[offset: offset + per_page]
    pagination_test_list_posts = get_test_list_posts(offset=offset, per_page=per_page)
    pagination = Pagination(page=page, per_page=per_page, total=total,
                        css_framework='bootstrap4',
                        display_msg="Показано <b>{start} - {end}</b> {record_name} из <b>{total}</b>",
                        record_name="записей")
    return flask.render_template("test/test_list_posts.html",
                                 test_list_posts=pagination_test_list_posts, 
                                                            page=page,
                           per_page=per_page,
                           pagination=pagination,)
@app.route("/view/<int:test_id>")
def get_post_test_post(test_id):
    conn = get_db_connection()
    with conn.cursor() as cur:
        sql = "SELECT * FROM `test_table` WHERE `test_id` =%s"
        cur.execute(sql, test_id)
        test_view_post = cur.fetchone()
    conn.close()
    return flask.render_template("test/test_view_post.html",
                                 test_view_post=test_view_post, )
@app.route("/edit/<int:test_id>/", methods=("GET", "POST"))
def edit_test_post(test_id):
    conn = get_db_connection()
    with conn.cursor() as cur:
        sql = "SELECT * FROM `test_table` WHERE `test_id` =%s"
        cur.execute(sql, (test_id,))
        edit_test_view = cur.fetchone()
    if flask.request.method == "POST":
        test_edit_post_text = flask.request.form["test_text"]
        if len(flask.request.form['test_text']) > 1:
            conn = get_db_connection()
            with conn.cursor() as cur:
                sql = "UPDATE `test_table` SET `test_text` =%s  WHERE `test_id` =%s"
                cur.execute(
                    sql, (test_edit_post_text, test_id),
                )
            conn.commit()
            conn.close()
            if not test_edit_post_text:
                flask.flash('Ошибка сохранения записи, вы ввели мало символов!', category='error')
            else:
                flask.flash('Запись успешно сохранена!', category='success')
            # В случае соблюдения условий заполнения полей, произойдёт перенаправление
            return flask.redirect(flask.url_for("index"))
        else:
            flask.flash('Ошибка сохранения записи!', category='error')
    return flask.render_template("test/edit_test_post.html", edit_test_view=edit_test_view)
@app.route("/new_post", methods=["GET", "POST"])
def add_test_post():
    if flask.request.method == "POST":
        new_test_post = flask.request.form["test_text"]
        if len(flask.request.form['test_text']) > 1:
            conn = get_db_connection()
            with conn.cursor() as cur:
                sql = "INSERT INTO `test_table` (`test_text`) VALUES (%s)"
                cur.execute(
                    sql, new_test_post,
                )
            conn.commit()
            conn.close()
            if not new_test_post:
                flask.flash('Ошибка сохранения записи, Вы ввели слишком мало символов!', category='error')
            else:
                flask.flash('Запись успешно добавлена!')
            # В случае соблюдения условий заполнения полей, произойдёт перенаправление
            return flask.redirect(flask.url_for("index"))
        else:
            flask.flash('Ошибка сохранения записи!', category='error')
    return flask.render_template("test/add_test_post.html")
@app.route("/delete/<int:test_id>/", methods=("POST",))
def delete_post_test(test_id):
    conn = get_db_connection()
    with conn.cursor() as cur:
        sql = "DELETE FROM `test_table` WHERE `test_id` =%s"
        cur.execute(
            sql, test_id,
        )
    conn.commit()
    conn.close()
    return flask.redirect(flask.url_for("index"))
if __name__ == "__main__":
    app.run(debug=True, host='0.0.0.0', port=83)
'@

# ---- Sheet: Bash (new rows 67-70) ----
$wsBash = $wb.Worksheets.Item("Bash")

$wsBash.Range("A67").Value = 97
$wsBash.Range("B67").Value = @'
#!/bin/commands
rm -rf ~/Downloads
mkdir ~/Downloads
'@
$wsBash.Range("C67").Value = @'
Скрипт для удаления папки и её создания

'@

$wsBash.Range("A68").Value = 98
$wsBash.Range("B68").Value = @'
#!/bin/bash 
# Задаем переменную DATE, в переменной указываем команду date (для вывода даты и времени) 
DATE=$(date) 
# Задаем переменную TEXT, в переменной указываем текст коммита и добавляем ранее добавленную переменную 
DATE TEXT="Add commit date: $DATE" 
# Добавляем файлы для коммита 
git add . 
# Объявляем коммит с созданной ранее переменной TEXT 
git commit -m "$TEXT" 
# Отправляем коммит на репозиторий 
git push 
'@
$wsBash.Range("C68").Value = @'
bash-скрипт для формирования коммита и отправки его на репозиторий с указанием даты 

'@

$wsBash.Range("A69").Value = 99
$wsBash.Range("B69").Value = @'
#!/bin/bash
# Запускается 1-е приложение следующим скриптом. Указывается путь исполняемой программы и путь к самому приложению.
~/AppData/Local/Programs/Python/Python312/python.exe ~/flask-mysql-project/app.py &
# Затем запускается 2-е приложение
~/AppData/Local/Programs/Python/Python312/python.exe /p/s.savin/flask-base/app.py &
# И наконец запускается 3-е приложение
~/AppData/Local/Programs/Python/Python312/python.exe /p/s.savin/flask-project-full/app.py &
'@
$wsBash.Range("C69").Value = @'
bash-скрипт который запускает три приложения flask

'@

$wsBash.Range("A70").Value = 100
$wsBash.Range("B70").Value = @'
sfsdfdsfsdfsdfsd
'@
$wsBash.Range("C70").Value = @'
sdfsdfsdfsdfsdfsdfsd
'@
